$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Populate new rows 19-42 with test case data (values first).
$ws.Range("A19").Value = "JinZu-ApiEngine-Test-7-var1"
$ws.Range("B19").Value = "good request, data retrieved"
$ws.Range("C19").Value = "{Project(cond:`"{charge_frequency:{_eq:3}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D19").Value = 200
$ws.Range("E19").Value = 100000
$ws.Range("F19").Value = "Successfully"

$ws.Range("A20").Value = "JinZu-ApiEngine-Test-7-var2"
$ws.Range("B20").Value = "good request, data retrieved"
$ws.Range("C20").Value = "{Project(cond:`"{is_manufacture_buy_back:{_eq:false}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D20").Value = 200
$ws.Range("E20").Value = 100000
$ws.Range("F20").Value = "Successfully"

$ws.Range("A21").Value = "JinZu-ApiEngine-Test-7-var4"
$ws.Range("B21").Value = "good request, data retrieved"
$ws.Range("C21").Value = "{Project(cond:`"{business_mgr:{_eq:\`"胡晓峰\`"}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D21").Value = 200
$ws.Range("E21").Value = 100000
$ws.Range("F21").Value = "Successfully"

$ws.Range("A22").Value = "JinZu-ApiEngine-Test-7-var5"
$ws.Range("B22").Value = "good request, data retrieved"
$ws.Range("C22").Value = "{Project(cond:`"{charge_frequency:{_neq:3}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type}}"
$ws.Range("D22").Value = 200
$ws.Range("E22").Value = 100000
$ws.Range("F22").Value = "Successfully"

$ws.Range("A23").Value = "JinZu-ApiEngine-Test-7-var6"
$ws.Range("B23").Value = "good request, data retrieved"
$ws.Range("C23").Value = "{Project(cond:`"{is_manufacture_buy_back:{_neq:false}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D23").Value = 200
$ws.Range("E23").Value = 100000
$ws.Range("F23").Value = "Successfully"

$ws.Range("A24").Value = "JinZu-ApiEngine-Test-7-var7"
$ws.Range("B24").Value = "good request, data retrieved"
$ws.Range("C24").Value = "{Project(cond:`"{status:{_neq:\`"online\`"}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D24").Value = 200
$ws.Range("E24").Value = 100000
$ws.Range("F24").Value = "Successfully"

$ws.Range("A25").Value = "JinZu-ApiEngine-Test-7-var8"
$ws.Range("B25").Value = "good request, data retrieved"
$ws.Range("C25").Value = "{Project(cond:`"{business_mgr:{_neq:\`"胡晓峰\`"}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D25").Value = 200
$ws.Range("E25").Value = 100000
$ws.Range("F25").Value = "Successfully"

$ws.Range("A26").Value = "JinZu-ApiEngine-Test-7-var9"
$ws.Range("B26").Value = "good request, data retrieved"
$ws.Range("C26").Value = "{Project(cond:`"{charge_frequency:{_gte:3}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D26").Value = 200
$ws.Range("E26").Value = 100000
$ws.Range("F26").Value = "Successfully"

$ws.Range("A27").Value = "JinZu-ApiEngine-Test-7-var10"
$ws.Range("B27").Value = "good request, data retrieved"
$ws.Range("C27").Value = "{Project(cond:`"{status: {_gte: \`"b\`"}}`") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}"
$ws.Range("D27").Value = 200
$ws.Range("E27").Value = 100000
$ws.Range("F27").Value = "Successfully"

$ws.Range("A28").Value = "JinZu-ApiEngine-Test-7-var11"
$ws.Range("B28").Value = "good request, data retrieved"
$ws.Range("C28").Value = "{Project(cond:`"{charge_frequency:{_gt:3}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D28").Value = 200
$ws.Range("E28").Value = 100000
$ws.Range("F28").Value = "Successfully"

$ws.Range("A29").Value = "JinZu-ApiEngine-Test-7-var12"
$ws.Range("B29").Value = "good request, data retrieved"
$ws.Range("C29").Value = "{Project(cond:`"{status: {_gt: \`"b\`"}}`") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}"
$ws.Range("D29").Value = 200
$ws.Range("E29").Value = 100000
$ws.Range("F29").Value = "Successfully"

$ws.Range("A30").Value = "JinZu-ApiEngine-Test-7-var13"
$ws.Range("B30").Value = "good request, data retrieved"
$ws.Range("C30").Value = "{Project(cond:`"{charge_frequency:{_lte:3}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D30").Value = 200
$ws.Range("E30").Value = 100000
$ws.Range("F30").Value = "Successfully"

$ws.Range("A31").Value = "JinZu-ApiEngine-Test-7-var14"
$ws.Range("B31").Value = "good request, data retrieved"
$ws.Range("C31").Value = "{Project(cond:`"{status: {_lte: \`"b\`"}}`") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}"
$ws.Range("D31").Value = 200
$ws.Range("E31").Value = 100000
$ws.Range("F31").Value = "Successfully"

$ws.Range("A32").Value = "JinZu-ApiEngine-Test-7-var15"
$ws.Range("B32").Value = "good request, data retrieved"
$ws.Range("C32").Value = "{Project(cond:`"{charge_frequency:{_lt:3}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D32").Value = 200
$ws.Range("E32").Value = 100000
$ws.Range("F32").Value = "Successfully"

$ws.Range("A33").Value = "JinZu-ApiEngine-Test-7-var16"
$ws.Range("B33").Value = "good request, data retrieved"
$ws.Range("C33").Value = "{Project(cond:`"{status: {_lt: \`"b\`"}}`") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}"
$ws.Range("D33").Value = 200
$ws.Range("E33").Value = 100000
$ws.Range("F33").Value = "Successfully"

$ws.Range("A34").Value = "JinZu-ApiEngine-Test-7-var17"
$ws.Range("B34").Value = "good request, data retrieved"
$ws.Range("C34").Value = "{Project (cond:`"{charge_frequency:{_in:[1,2,3]}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D34").Value = 200
$ws.Range("E34").Value = 100000
$ws.Range("F34").Value = "Successfully"

$ws.Range("A35").Value = "JinZu-ApiEngine-Test-7-var18"
$ws.Range("B35").Value = "good request, data retrieved"
$ws.Range("C35").Value = "{Project(cond:`"{status:{_in:[\`"online\`",\`"archived\`"]}}`",order:`"`") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}"
$ws.Range("D35").Value = 200
$ws.Range("E35").Value = 100000
$ws.Range("F35").Value = "Successfully"

$ws.Range("A36").Value = "JinZu-ApiEngine-Test-7-var19"
$ws.Range("B36").Value = "good request, data retrieved"
$ws.Range("C36").Value = "{Project(cond:`"{business_mgr:{_in:[\`"潘云晖\`",\`"臧佳宝\`"]}}`",order:`"`") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}"
$ws.Range("D36").Value = 200
$ws.Range("E36").Value = 100000
$ws.Range("F36").Value = "Successfully"

$ws.Range("A37").Value = "JinZu-ApiEngine-Test-7-var20"
$ws.Range("B37").Value = "good request, data retrieved"
$ws.Range("C37").Value = "{Project (cond:`"{charge_frequency:{_nin:[1,2,3]}}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D37").Value = 200
$ws.Range("E37").Value = 100000
$ws.Range("F37").Value = "Successfully"

$ws.Range("A38").Value = "JinZu-ApiEngine-Test-7-var21"
$ws.Range("B38").Value = "good request, data retrieved"
$ws.Range("C38").Value = "{Project(cond:`"{status:{_nin:[\`"online\`",\`"bbbb\`"]}}`",order:`"`") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}"
$ws.Range("D38").Value = 200
$ws.Range("E38").Value = 100000
$ws.Range("F38").Value = "Successfully"

$ws.Range("A39").Value = "JinZu-ApiEngine-Test-7-var22"
$ws.Range("B39").Value = "good request, data retrieved"
$ws.Range("C39").Value = "{Project(cond:`"{business_mgr:{_nin:[\`"潘云晖\`",\`"臧佳宝\`"]}}`",order:`"`") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}"
$ws.Range("D39").Value = 200
$ws.Range("E39").Value = 100000
$ws.Range("F39").Value = "Successfully"

$ws.Range("A40").Value = "JinZu-ApiEngine-Test-7-var23"
$ws.Range("B40").Value = "good request, data retrieved"
$ws.Range("C40").Value = "{Project (cond:`"{_and: [{business_mgr:{_in:[\`"潘云晖\`",\`"臧佳宝\`" ]}},{status:{_in:[\`"archived\`" ]}}]}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D40").Value = 200
$ws.Range("E40").Value = 100000
$ws.Range("F40").Value = "Successfully"

$ws.Range("A41").Value = "JinZu-ApiEngine-Test-7-var24"
$ws.Range("B41").Value = "good request, data retrieved"
$ws.Range("C41").Value = "{Project (cond:`"{_or: [{business_mgr:{_in:[\`"潘云晖\`",\`"臧佳宝\`"]}},{status:{_in:[\`"archived\`" ]}}]}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}"
$ws.Range("D41").Value = 200
$ws.Range("E41").Value = 100000
$ws.Range("F41").Value = "Successfully"

$ws.Range("A42").Value = "JinZu-ApiEngine-Test-7-var25"
$ws.Range("B42").Value = "good request, data retrieved"
$ws.Range("C42").Value = "{Project(cond:`"{_or: [{business_mgr:{_in:[\`"潘云晖\`",\`"臧佳宝\`" ]}},{status:{_in:[\`"archived\`" ]}}]}`",order:`"`") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status  invert_Customer (cond:`"{id:{_eq:24}}`",order:`"`") {actual_controller category cid city cname contact contact_detail ctype district enterprise_size group holding_type id is_connected_tx is_gov_fin_customer is_group_customer legal_person legal_person_id major_class middle_class office_address project province registered_address small_class}}}"
$ws.Range("D42").Value = 200
$ws.Range("E42").Value = 100000
$ws.Range("F42").Value = "Successfully"

# Copy formatting from the last pre-existing row (18) down onto the new rows
# so the new cells inherit the same font/border/fill as the rest of the table
# (mirrors how the row was originally authored by typing below row 18 in Excel).
$ws.Range("A18:F18").Copy()
$ws.Range("A19:F42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column C needed to widen to fit the new, longer GraphQL query strings.
$ws.Columns.Item(3).ColumnWidth = 90.7109375

# Restore the selection state recorded for this sheet after the edit.
$ws.Range("C33").Select()
